$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings: week number + date range) ---
$ws.Range("A8").Value = "Volume 30   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Precinct crime-statistics table updates (rows 15-21, 24-30) ---
$ws.Range("L15").Value = -54.166666666666
$ws.Range("M15").Value = -8.333333333333
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -60
$ws.Range("I16").Value = 93
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = -16.216216216216
$ws.Range("L16").Value = 20.779220779220
$ws.Range("M16").Value = -51.813471502590
$ws.Range("N16").Value = -84.048027444253
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 256
$ws.Range("J17").Value = 263
$ws.Range("K17").Value = -2.661596958174
$ws.Range("L17").Value = 6.224066390041
$ws.Range("M17").Value = 43.016759776536
$ws.Range("N17").Value = 12.280701754386
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = -28.571428571428
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 13.333333333333
$ws.Range("I18").Value = 120
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 30.434782608695
$ws.Range("L18").Value = 48.148148148148
$ws.Range("M18").Value = -42.028985507246
$ws.Range("N18").Value = -86.984815618221
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 30.612244897959
$ws.Range("I19").Value = 373
$ws.Range("J19").Value = 369
$ws.Range("K19").Value = 1.084010840108
$ws.Range("L19").Value = 51.626016260162
$ws.Range("M19").Value = 45.703125
$ws.Range("N19").Value = 6.876790830945
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -4.347826086956
$ws.Range("I20").Value = 186
$ws.Range("J20").Value = 133
$ws.Range("K20").Value = 39.849624060150
$ws.Range("L20").Value = 111.363636363636
$ws.Range("M20").Value = -15.068493150684
$ws.Range("N20").Value = -90.332640332640
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 21.875
$ws.Range("F21").Value = 149
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = -2.614379084967
$ws.Range("I21").Value = 1041
$ws.Range("J21").Value = 978
$ws.Range("K21").Value = 6.441717791411
$ws.Range("L21").Value = 36.793692509855
$ws.Range("M21").Value = -3.162790697674
$ws.Range("N21").Value = -74.245423057892
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = 8.571428571428
$ws.Range("I24").Value = 711
$ws.Range("J24").Value = 849
$ws.Range("K24").Value = -16.254416961130
$ws.Range("L24").Value = 30.458715596330
$ws.Range("M24").Value = 48.125
$ws.Range("C25").Value = 10
$ws.Range("E25").Value = -23.076923076923
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = -1.785714285714
$ws.Range("I25").Value = 364
$ws.Range("J25").Value = 391
$ws.Range("K25").Value = -6.905370843989
$ws.Range("L25").Value = 41.085271317829
$ws.Range("M25").Value = -5.208333333333
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = -34.375
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -62.5
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -9.677419354838
$ws.Range("L27").Value = -3.448275862068
$ws.Range("G28").Value = 6
$ws.Range("N28").Value = -81.578947368421
$ws.Range("G29").Value = 4
$ws.Range("N29").Value = -83.783783783783
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = -50
